$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for this market/category. In this
# worksheet each record occupies one row, and existing rows 157-193 need to
# shift "down" by one row to make room for the brand new record which lands
# in row 157, pushing the former row 193 down into a brand new row 194.

# 1) Clone row 193 (last row of the block) into the new row 194 so that all
#    the columns that stay constant across the block (A, B, C, E, F, G, H,
#    I, L, N, O, Q, R) - as well as D/J/K/M/P, which will correctly hold the
#    former row 193 values - are populated and correctly styled.
$ws.Range("A193:R193").Copy($ws.Range("A194:R194"))

# 2) Shift the date/price columns (D, J, K, M, P) down by one row for rows
#    158..193, each one taking the value that used to belong to the row
#    directly above it. Processing from the bottom (193) upward to the top
#    (158) guarantees the source row for each step has not yet been
#    overwritten.
for ($i = 193; $i -ge 158; $i--) {
    $src = $i - 1
    $ws.Cells.Item($i, 4).Value = $ws.Cells.Item($src, 4).Value2()
    $ws.Cells.Item($i, 10).Value = $ws.Cells.Item($src, 10).Value2()
    $ws.Cells.Item($i, 11).Value = $ws.Cells.Item($src, 11).Value2()
    $ws.Cells.Item($i, 13).Value = $ws.Cells.Item($src, 13).Value2()
    $ws.Cells.Item($i, 16).Value = $ws.Cells.Item($src, 16).Value2()
    $ws.Cells.Item($i, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# 3) Row 157 now holds the brand new record's values.
$ws.Cells.Item(157, 4).Value = 44543
$ws.Cells.Item(157, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(157, 10).Value = 2300
$ws.Cells.Item(157, 11).Value = 400
$ws.Cells.Item(157, 13).Value = 450
$ws.Cells.Item(157, 16).Value = 900
